$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text values - Excel won't reinterpret these as numbers/dates.
$ws.Range("A1").Value = "None"
$ws.Range("B1").Value = "None"
$ws.Range("C1").Value = "None"

# "0.0" and "09/08/2024" look like a number / a date to Excel, so force
# the cell to Text format before typing them in, then drop the
# formatting again so the cell keeps using the workbook's default style.
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "0.0"
$ws.Range("E1").ClearFormats()

$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Value = "09/08/2024"
$ws.Range("H1").ClearFormats()

$ws.Range("G1").Value = "Cigna"
